$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the whole sheet (content + formatting) so stale cells (old column F,
# old column E data, etc.) do not linger as empty-but-styled cells.
$ws.Cells.Clear()

# ---- Header row ----
$ws.Cells.Item(1,1).Value = "Parameter"
$ws.Cells.Item(1,2).Value = "Subparameter"
$ws.Cells.Item(1,3).Value = "Value"
$ws.Cells.Item(1,4).Value = "Enabled"
$ws.Cells.Item(1,5).Value = "Default"

# ---- Data rows ----
# Row 2: growth / threemonths - no Value, Enabled = FALSE; Value cell is
# formatted as a percentage but left blank.
$ws.Cells.Item(2,1).Value = "growth"
$ws.Cells.Item(2,2).Value = "threemonths"
$ws.Cells.Item(2,4).Value = $false
$ws.Cells.Item(2,3).NumberFormat = "0%"

# Row 3: growth / oneyear - no Value, Enabled = FALSE
$ws.Cells.Item(3,1).Value = "growth"
$ws.Cells.Item(3,2).Value = "oneyear"
$ws.Cells.Item(3,4).Value = $false

# Row 4: sma10 / length, Value = 10, Enabled = TRUE
$ws.Cells.Item(4,1).Value = "sma10"
$ws.Cells.Item(4,2).Value = "length"
$ws.Cells.Item(4,3).Value = 10
$ws.Cells.Item(4,4).Value = $true

# Row 5: stoch / k, Value = 14, Enabled = TRUE
$ws.Cells.Item(5,1).Value = "stoch"
$ws.Cells.Item(5,2).Value = "k"
$ws.Cells.Item(5,3).Value = 14
$ws.Cells.Item(5,4).Value = $true

# Row 6: stoch / d, Value = 3, Enabled = TRUE
$ws.Cells.Item(6,1).Value = "stoch"
$ws.Cells.Item(6,2).Value = "d"
$ws.Cells.Item(6,3).Value = 3
$ws.Cells.Item(6,4).Value = $true

# Row 7: stoch / klen, Value = 3, Enabled = TRUE
$ws.Cells.Item(7,1).Value = "stoch"
$ws.Cells.Item(7,2).Value = "klen"
$ws.Cells.Item(7,3).Value = 3
$ws.Cells.Item(7,4).Value = $true

# Row 8: rsi / length, Value = 14, Enabled = TRUE
$ws.Cells.Item(8,1).Value = "rsi"
$ws.Cells.Item(8,2).Value = "length"
$ws.Cells.Item(8,3).Value = 14
$ws.Cells.Item(8,4).Value = $true

# Row 9: atrts / k, Value = 3, Enabled = TRUE
$ws.Cells.Item(9,1).Value = "atrts"
$ws.Cells.Item(9,2).Value = "k"
$ws.Cells.Item(9,3).Value = 3
$ws.Cells.Item(9,4).Value = $true

# Row 10: atrts / length, Value = 14, Enabled = TRUE
$ws.Cells.Item(10,1).Value = "atrts"
$ws.Cells.Item(10,2).Value = "length"
$ws.Cells.Item(10,3).Value = 14
$ws.Cells.Item(10,4).Value = $true

# Row 11: atr / length, Value = 14, Enabled = TRUE
$ws.Cells.Item(11,1).Value = "atr"
$ws.Cells.Item(11,2).Value = "length"
$ws.Cells.Item(11,3).Value = 14
$ws.Cells.Item(11,4).Value = $true

# ---- Selection state ----
# Target workbook ends up with a multi-area selection (E3, F11) with F11 as
# the active cell; this runtime only models a single active range, so get as
# close as possible by leaving the active cell on F11.
$ws.Range("E3").Select()
$ws.Range("F11").Select()
